# NatureTool Input File Template - "Added][Documentation] Added a Readme.md file,
# updated the path to results file and made a first run"
#
# This script fills in the previously-empty "HELP" sheet with a documentation
# block describing the Input_data sheet, tweaks one ISIN value (and its
# formatting) on the Input_data sheet, and restores the active-cell selection
# on both sheets.

$wb = $excel.ActiveWorkbook

$help = $wb.Worksheets.Item("HELP")
$data = $wb.Worksheets.Item("Input_data")

# ---------------------------------------------------------------------------
# HELP sheet - title banner (B1:C1)
# ---------------------------------------------------------------------------
$help.Range("B1").Value = "NATURE TOOL INPUT FILE"

$titleBand = $help.Range("B1:C1")
$titleBand.Interior.ThemeColor = 3      # -> OOXML theme="3" (Text 2 / dk2)
$titleBand.Font.Bold = $true
$titleBand.Font.ThemeColor = 2          # -> OOXML theme="0" (Background 1 / lt1, white)

# ---------------------------------------------------------------------------
# HELP sheet - intro sentence (B2), partially italicised
# ---------------------------------------------------------------------------
$help.Range("B2").Value = 'Please fill in the "Input_data" sheet before running the NatureTool'
$introItalic = $help.Range("B2").Characters(20, 13)
$introItalic.Font.Italic = $true
$introRest = $help.Range("B2").Characters(33, 36)
$introRest.Font.Italic = $false

# ---------------------------------------------------------------------------
# HELP sheet - "Input_data" sub-header (B4)
# ---------------------------------------------------------------------------
$help.Range("B4").Value = "Input_data"
$help.Range("B4").Font.Bold = $true
$help.Range("B4").Interior.ThemeColor = 8   # -> OOXML theme="7" (Accent 4), same as the Input_data tab color

# ---------------------------------------------------------------------------
# HELP sheet - column glossary (rows 5-8)
# ---------------------------------------------------------------------------
$help.Range("B5").Value = "Portfolio_name"
$help.Range("B5").Font.Bold = $true
$help.Range("C5").Value = "Name of your portfolio (can be the same for all lines)"

$help.Range("B6").Value = "ISIN"
$help.Range("B6").Font.Bold = $true
$help.Range("C6").Value = "ISIN number in the portfolio"

$help.Range("B7").Value = "Amount"
$help.Range("B7").Font.Bold = $true
$help.Range("C7").Value = "Invested amount related to the ISIN (equity or bond)"
$help.Range("C7").Font.Size = 10

$help.Range("B8").Value = "Currency"
$help.Range("B8").Font.Bold = $true
$help.Range("C8").Value = 'Currency of the amount. The list of authorized currencies is defined in the "Currency" sheet'
$currencyItalic = $help.Range("C8").Characters(77, 10)
$currencyItalic.Font.Italic = $true
$currencyRest = $help.Range("C8").Characters(87, 6)
$currencyRest.Font.Italic = $false

# ---------------------------------------------------------------------------
# HELP sheet - column widths
# ---------------------------------------------------------------------------
$help.Columns.Item(1).ColumnWidth = 4.21875
$help.Columns.Item(2).ColumnWidth = 14.21875

# ---------------------------------------------------------------------------
# Input_data sheet - correct one ISIN and frame it with a thin border
# ---------------------------------------------------------------------------
$data.Range("B5").Value = "GB00BDCPN049"
$data.Range("B5").Borders.LineStyle = 1
$data.Range("B5").Borders.Weight = 2

# ---------------------------------------------------------------------------
# Restore selections: HELP -> C10, Input_data -> E6 (Input_data stays the
# active/visible tab, matching tabSelected="1" in the workbook)
# ---------------------------------------------------------------------------
[void]$help.Activate()
[void]$help.Range("C10").Select()

[void]$data.Activate()
[void]$data.Range("E6").Select()
